# Apply price / volume / coin-rank updates scraped on Thu Apr 18 12:43:23 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.931.48"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "3.008.88"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'540.85"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "'132.59"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "3.004.38"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "'0.491"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'6.10"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.146"
$ws.Range("E11").Value = "  -3.51%  "
$ws.Range("D12").Value = "'0.444"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "'0.0000220"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "'34.19"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "3.493.76"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "61.880.76"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.108"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.008.13"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "'6.62"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "'477.80"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "'13.19"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").Value = "'0.670"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "'7.02"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").Value = "'80.67"
$ws.Range("E24").Value = "  +4.62%  "
$ws.Range("D25").Value = "'12.06"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "'2.70"
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("D28").Value = "'7.69"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'1.92"
$ws.Range("E30").Value = "  +4.81%  "
$ws.Range("D31").Value = "'25.58"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "'1.12"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "'5.64"
$ws.Range("E33").Value = "  +6.53%  "
$ws.Range("D34").Value = "'2.33"
$ws.Range("E34").Value = "  +3.77%  "
$ws.Range("D35").Value = "'54.80"
$ws.Range("E35").Value = "  -5.75%  "
$ws.Range("D36").Value = "'5.84"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'449.37"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("D38").Value = "3.152.34"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").Value = "'0.0796"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").Value = "'0.0382"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "'8.06"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "'2.43"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").Value = "'26.23"
$ws.Range("E44").Value = "  +7.23%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'0.242"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "'0.108"
$ws.Range("D48").Value = "'1.95"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("D49").Value = "'114.31"
$ws.Range("E49").Value = "  -5.77%  "
$ws.Range("D50").Value = "0.0₃0494"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("E51").Value = "  +4.09%  "
